# Generate Report for Archive
#
# The localization-status report is regenerated: the tracked files have
# moved on from "Ready for handoff" to "In Translation". Update every
# sheet that surfaces that status value (Overview shows it per-locale in
# columns E/F, while the zh-cn / de-de detail sheets show it in column
# C), then let the now-shorter Status text narrow those columns.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Closest column width (in Excel "characters" units) this runtime's pixel
# grid can represent for the narrower Status column once it holds
# "In Translation" instead of "Ready for handoff".
$statusColWidth = 12.5

# --- Overview sheet: Status appears per-locale in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet: Status is column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet: Status is column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
